# Fruta / hortaliza, semanal
# Insert a new weekly record as row 258 (pushing the existing rows 258-281
# down to 259-282) on the single data sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 258, shifting all the
# following rows (old 258..281) down to (259..282).
$ws.Rows.Item(258).Insert()

# Populate the newly inserted row 258 with the new weekly observation.
$ws.Cells.Item(258, 1).Value = 11
$ws.Cells.Item(258, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(258, 3).Value = "Bíobío"
$ws.Cells.Item(258, 4).Value = 45212
$ws.Cells.Item(258, 5).Value = 8
$ws.Cells.Item(258, 6).Value = 100112032
$ws.Cells.Item(258, 7).Value = "Zapallo italiano"
$ws.Cells.Item(258, 8).Value = "Sin especificar"
$ws.Cells.Item(258, 9).Value = "Primera"
$ws.Cells.Item(258, 10).Value = 270
$ws.Cells.Item(258, 11).Value = 18000
$ws.Cells.Item(258, 12).Value = 20000
$ws.Cells.Item(258, 13).Value = 19111
$ws.Cells.Item(258, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(258, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(258, 16).Value = 382
$ws.Cells.Item(258, 17).Value = 50
$ws.Cells.Item(258, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Cells.Item(258, 4).NumberFormat = $ws.Cells.Item(259, 4).NumberFormat()
